$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("Q2").Value = 3.4
$ws.Range("R2").Value = 1.33

# Row 3 updates
$ws.Range("I3").Value = 5.75
$ws.Range("J3").Value = 1.95
$ws.Range("K3").Value = 2.6
$ws.Range("Q3").Value = 1.48
$ws.Range("R3").Value = 2.6
$ws.Range("U3").Value = 1.67
$ws.Range("V3").Value = 2.1
$ws.Range("AD3").Value = 10
$ws.Range("AE3").Value = 17
$ws.Range("AQ3").Value = 19
$ws.Range("AZ3").Value = 101
